$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.238.57'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.04%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.664.35'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.03%  '

$ws.Range('E4').Value = '  -0.44%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '219.92'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.26%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5263'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.17%  '

$ws.Range('E7').Value = '  -0.40%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2684'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.75%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06395'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.17%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.78'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.18%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07704'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.73%  '

$ws.Range('B12').NumberFormat = '@'
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('B12').Style = 'Normal'
$ws.Range('C12').NumberFormat = '@'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('C12').Style = 'Normal'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.670'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.85%  '

$ws.Range('B13').NumberFormat = '@'
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('B13').Style = 'Normal'
$ws.Range('C13').NumberFormat = '@'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('C13').Style = 'Normal'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.648.42'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.96%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.893.13'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.07%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5664'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.88%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0₅8281'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.08%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '65.84'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.95%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '26.241.54'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.07%  '

$ws.Range('E19').Value = '  -0.43%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.706'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.82%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.61'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.23%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '193.60'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.04%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.017'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.64%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.004'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.49%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.36'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.35%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1208'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.28%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.341'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.26%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '16.14'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.19%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.523'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.20%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05639'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.53%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.281'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.14%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.504'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.53%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.402'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.17%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.577'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.45%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9617'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.60%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.787'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.91%  '

$ws.Range('E37').Value = '  -1.07%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5770'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.13%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01607'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.59%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.959'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.07%  '

$ws.Range('E41').Value = '  -0.42%  '

$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'Maker'
$ws.Range('B42').Style = 'Normal'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('C42').Style = 'Normal'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.041.62'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.72%  '

$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('B43').Style = 'Normal'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('C43').Style = 'Normal'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8393'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.32%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.59'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.45%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.803.77'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.08%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '58.80'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.57%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0₈104'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.71%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.004'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.63%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.126'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.67%  '

$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'Mantle'
$ws.Range('B50').Style = 'Normal'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('C50').Style = 'Normal'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4347'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.48%  '

$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'Cronos'
$ws.Range('B51').Style = 'Normal'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('C51').Style = 'Normal'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05245'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.65%  '
